$wb = $excel.ActiveWorkbook

# "Other high-grade glioma" sheet: update column C (pvalue) rows 2-9
$wsHighGrade = $wb.Worksheets.Item("Other high-grade glioma")
$wsHighGrade.Range("C2").Value = 0.0874616484567612
$wsHighGrade.Range("C3").Value = 0.0722971792396567
$wsHighGrade.Range("C4").Value = 0.632646633329443
$wsHighGrade.Range("C5").Value = 0.747469115830449
$wsHighGrade.Range("C6").Value = 0.765044002655438
$wsHighGrade.Range("C7").Value = 0.688070088451452
$wsHighGrade.Range("C8").Value = 0.806157348933911
$wsHighGrade.Range("C9").Value = 0.831034212389158

# "DIPG or DMG" sheet: update column C (pvalue) rows 2-9 (row 5 unchanged)
$wsDipg = $wb.Worksheets.Item("DIPG or DMG")
$wsDipg.Range("C2").Value = 0.158302905892242
$wsDipg.Range("C3").Value = 0.99999999999995
$wsDipg.Range("C4").Value = 0.999999999999968
$wsDipg.Range("C6").Value = 0.862159087147157
$wsDipg.Range("C7").Value = 0.895123569422763
$wsDipg.Range("C8").Value = 0.338624338624339
$wsDipg.Range("C9").Value = 0.574228093777768
